$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the rows for years 2000-2009 (rows 2-11). This shifts the
# 2010-2019 rows (formerly rows 12-21) up to rows 2-11.
$ws.Range("A2:E11").EntireRow.Delete()
